$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.790.23'
$ws.Range('E2').Value = '  +1.32%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.59'
$ws.Range('E3').Value = '  +1.42%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.032'
$ws.Range('E4').Value = '  +0.43%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.94'
$ws.Range('E5').Value = '  +1.46%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.028'
$ws.Range('E6').Value = '  +0.33%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4403'
$ws.Range('E7').Value = '  +0.97%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3812'
$ws.Range('E8').Value = '  +2.37%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07453'
$ws.Range('E9').Value = '  +1.19%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8888'
$ws.Range('E10').Value = '  +1.81%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.69'
$ws.Range('E11').Value = '  +1.41%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.880.79'
$ws.Range('E12').Value = '  -0.13%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.555'
$ws.Range('E13').Value = '  +1.71%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.766'
$ws.Range('E14').Value = '  +1.15%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07188'
$ws.Range('E15').Value = '  +0.65%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '85.44'
$ws.Range('E16').Value = '  +3.43%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.035'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009133'
$ws.Range('E18').Value = '  +1.60%  '

# Row 19
$ws.Range('E19').Value = '  +0.34%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.56'
$ws.Range('E20').Value = '  +1.00%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.793.66'
$ws.Range('E21').Value = '  +1.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.316'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.29'
$ws.Range('E23').Value = '  +1.13%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.105.04'
$ws.Range('E24').Value = '  +0.64%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.026'
$ws.Range('E25').Value = '  +6.83%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.60'
$ws.Range('E26').Value = '  +1.15%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.80'
$ws.Range('E27').Value = '  +1.14%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.400'
$ws.Range('E28').Value = '  +3.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.999'
$ws.Range('E29').Value = '  +3.93%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.03'
$ws.Range('E30').Value = '  +1.80%  '

# Row 31
$ws.Range('E31').Value = '  -0.28%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7852'
$ws.Range('E32').Value = '  +3.25%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.223'
$ws.Range('E33').Value = '  +1.83%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.026'
$ws.Range('E34').Value = '  +5.57%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.591'
$ws.Range('E35').Value = '  +2.54%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.031'

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.146'
$ws.Range('E37').Value = '  +0.05%  '

# Row 38
$ws.Range('E38').Value = '  +1.15%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05328'
$ws.Range('E39').Value = '  +1.54%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5227'
$ws.Range('E40').Value = '  +1.17%  '

# Row 41
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.858'
$ws.Range('E41').Value = '  +2.63%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1690'
$ws.Range('E42').Value = '  +1.72%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.924'
$ws.Range('E43').Value = '  +5.61%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.911'
$ws.Range('E44').Value = '  +4.97%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '111.02'
$ws.Range('E45').Value = '  +2.11%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.80'
$ws.Range('E46').Value = '  +2.46%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06604'
$ws.Range('E47').Value = '  +4.70%  '

# Row 48
$ws.Range('E48').Value = '  +2.60%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.031'
$ws.Range('E49').Value = '  +0.42%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4743'
$ws.Range('E50').Value = '  +2.53%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.923'
$ws.Range('E51').Value = '  +1.42%  '
